$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Hyperlinks.Delete()
$ws.Range("C1:D1").EntireColumn.Delete()
[void]$ws.Range("C1:D1048576").Select()
